# coletas.xlsx bugfix: excel & get-by-date
# Adds two new "PB-02" / "PB-03" measurement blocks (mirroring the existing
# "PB-01" block) to the "excel" sheet, updates the "Hora Fim" sample value,
# and appends a new data row representing a second collected record.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Duplicate the "PB-01" header block (row 1 + row 2) into two new
#        blocks, copying cell-by-cell so no stray blank cells are created
#        inside the soon-to-be-merged ranges. ---

# PB-02 block -> columns P:T
$ws.Range("H1").Copy($ws.Range("P1"))
$ws.Range("P1").Value = "PB-02"
$ws.Range("H2").Copy($ws.Range("P2"))
$ws.Range("I2").Copy($ws.Range("Q2"))
$ws.Range("J2").Copy($ws.Range("R2"))
$ws.Range("K2").Copy($ws.Range("S2"))
$ws.Range("L2").Copy($ws.Range("T2"))
$ws.Range("P1:T1").Merge()

# PB-03 block -> columns U:Y
$ws.Range("H1").Copy($ws.Range("U1"))
$ws.Range("U1").Value = "PB-03"
$ws.Range("H2").Copy($ws.Range("U2"))
$ws.Range("I2").Copy($ws.Range("V2"))
$ws.Range("J2").Copy($ws.Range("W2"))
$ws.Range("K2").Copy($ws.Range("X2"))
$ws.Range("L2").Copy($ws.Range("Y2"))
$ws.Range("U1:Y1").Merge()

# --- 2. Update existing row 3 data ---
# "Hora Fim" sample value changes
$ws.Range("D3").Value = "09:47:59.460"

# New PB-02 sample values for row 3 (copy style cell-by-cell, then overwrite)
$ws.Range("H3").Copy($ws.Range("P3"))
$ws.Range("P3").Value = 1
$ws.Range("I3").Copy($ws.Range("Q3"))
$ws.Range("Q3").Value = 2
$ws.Range("J3").Copy($ws.Range("R3"))
$ws.Range("R3").Value = 3
$ws.Range("K3").Copy($ws.Range("S3"))
$ws.Range("S3").Value = 4
$ws.Range("L3").Copy($ws.Range("T3"))
$ws.Range("T3").Value = 5

# --- 3. Append a new data row (row 4) for a second collected record ---
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("B3").Copy($ws.Range("B4"))
$ws.Range("C3").Copy($ws.Range("C4"))
$ws.Range("D3").Copy($ws.Range("D4"))
$ws.Range("E3").Copy($ws.Range("E4"))
$ws.Range("F3").Copy($ws.Range("F4"))
$ws.Range("G3").Copy($ws.Range("G4"))
$ws.Range("H3").Copy($ws.Range("H4"))
$ws.Range("I3").Copy($ws.Range("I4"))

$ws.Range("C4").Value = "10:00"
$ws.Range("D4").Value = "11:32:42.884"
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 10
$ws.Range("I4").Value = 10

# --- 4. Column widths for the newly introduced columns (cosmetic, bestFit) ---
$ws.Range("P1").ColumnWidth = 24.416666666666664
$ws.Range("U1").ColumnWidth = 24.416666666666664
$ws.Range("Q1").ColumnWidth = 15.25
$ws.Range("R1").ColumnWidth = 15.75
$ws.Range("S1").ColumnWidth = 14.75
$ws.Range("T1").ColumnWidth = 34.08333333333333
$ws.Range("V1").ColumnWidth = 15.25
$ws.Range("W1").ColumnWidth = 15.75
$ws.Range("X1").ColumnWidth = 14.75
$ws.Range("Y1").ColumnWidth = 34.08333333333333
